$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 1260, shifting all
# subsequent data (old rows 1260-1344) down to rows 1263-1347.
$ws.Rows.Item(1260).Insert()
$ws.Rows.Item(1260).Insert()
$ws.Rows.Item(1260).Insert()

# Columns A, B, C, E, F, G, R are constant across every data row in this
# sheet, so copy them down into the 3 freshly inserted rows from the row
# directly above (row 1259, still the last of the original data block).
for ($r = 1260; $r -le 1262; $r++) {
    $ws.Cells.Item($r, 1).Value2  = $ws.Cells.Item(1259, 1).Value2   # A - Mercado ID
    $ws.Cells.Item($r, 2).Value2  = $ws.Cells.Item(1259, 2).Value2   # B - Mercado
    $ws.Cells.Item($r, 3).Value2  = $ws.Cells.Item(1259, 3).Value2   # C - Region
    $ws.Cells.Item($r, 5).Value2  = $ws.Cells.Item(1259, 5).Value2   # E - Codreg
    $ws.Cells.Item($r, 6).Value2  = $ws.Cells.Item(1259, 6).Value2   # F - Categoria ID
    $ws.Cells.Item($r, 7).Value2  = $ws.Cells.Item(1259, 7).Value2   # G - Categoria
    $ws.Cells.Item($r, 18).Value2 = $ws.Cells.Item(1259, 18).Value2  # R - Clasificacion
}

# Row 1260
$ws.Cells.Item(1260, 4).Value2  = 44931
$ws.Cells.Item(1260, 8).Value   = "Morada(o)"
$ws.Cells.Item(1260, 9).Value   = "1a (guarda)"
$ws.Cells.Item(1260, 10).Value2 = 200
$ws.Cells.Item(1260, 11).Value2 = 10000
$ws.Cells.Item(1260, 12).Value2 = 10000
$ws.Cells.Item(1260, 13).Value2 = 10000
$ws.Cells.Item(1260, 14).Value  = "$/malla 18 kilos"
$ws.Cells.Item(1260, 15).Value  = "Región de O'Higgins"
$ws.Cells.Item(1260, 16).Value2 = 556
$ws.Cells.Item(1260, 17).Value2 = 18

# Row 1261
$ws.Cells.Item(1261, 4).Value2  = 44931
$ws.Cells.Item(1261, 8).Value   = "Sin especificar"
$ws.Cells.Item(1261, 9).Value   = "1a nueva(o)"
$ws.Cells.Item(1261, 10).Value2 = 1900
$ws.Cells.Item(1261, 11).Value2 = 8500
$ws.Cells.Item(1261, 12).Value2 = 10000
$ws.Cells.Item(1261, 13).Value2 = 9211
$ws.Cells.Item(1261, 14).Value  = "$/malla 18 kilos"
$ws.Cells.Item(1261, 15).Value  = "Región del Maule"
$ws.Cells.Item(1261, 16).Value2 = 512
$ws.Cells.Item(1261, 17).Value2 = 18

# Row 1262
$ws.Cells.Item(1262, 4).Value2  = 44931
$ws.Cells.Item(1262, 8).Value   = "Sin especificar"
$ws.Cells.Item(1262, 9).Value   = "Primera"
$ws.Cells.Item(1262, 10).Value2 = 500
$ws.Cells.Item(1262, 11).Value2 = 10000
$ws.Cells.Item(1262, 12).Value2 = 10000
$ws.Cells.Item(1262, 13).Value2 = 10000
$ws.Cells.Item(1262, 14).Value  = "$/malla 18 kilos"
$ws.Cells.Item(1262, 15).Value  = "Perú"
$ws.Cells.Item(1262, 16).Value2 = 556
$ws.Cells.Item(1262, 17).Value2 = 18
